$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.420.37'
$ws.Range('D2').ClearFormats()
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.484.38'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.78%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.54'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.14'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.57%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.92%  '
$ws.Range('E9').Value = '  +4.61%  '
$ws.Range('E10').Value = '  +0.52%  '
$ws.Range('E11').Value = '  +3.99%  '
$ws.Range('E12').Value = '  +2.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.31'
$ws.Range('D13').ClearFormats()
$ws.Range('E14').Value = '  +6.68%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.948.80'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +3.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.274.99'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.485.62'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.59'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.26'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +6.57%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +2.51%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '328.04'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('E23').Value = '  +10.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.61'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '637.49'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +16.01%  '
$ws.Range('E26').Value = '  +13.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.76'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('E28').Value = '  +4.75%  '
$ws.Range('E29').Value = '  +9.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.44'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.54%  '
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -1.92%  '
$ws.Range('E33').Value = '  +2.80%  '
$ws.Range('E34').Value = '  +10.38%  '
$ws.Range('E35').Value = '  +3.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('E37').Value = '  +2.13%  '
$ws.Range('E38').Value = '  +1.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.94'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.39%  '
$ws.Range('E40').Value = '  +3.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '147.05'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.14%  '
$ws.Range('E42').Value = '  +17.90%  '
$ws.Range('E43').Value = '  +0.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '150.53'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.75%  '
$ws.Range('E45').Value = '  +3.89%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.18'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +6.80%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0549'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +4.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.608'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0239'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +5.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0928'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('E51').Value = '  +5.59%  '
